$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CharacterAtk table data updates (row 3 of the CharacterAtkTable)
# ThrowCooldown 2 -> 20
$ws.Range("D3").Value = 20
# SwingCooldown 1 -> 20
$ws.Range("H3").Value = 20
# SwingRad 2.5 -> 1.5
$ws.Range("J3").Value = 1.5

# Update the active selection shown in the saved view
$ws.Range("E10").Select()
